# Append the new trade (Trade #34, closed/opened 2026-02-16 22:55:29,
# base_strategy DOWN) as row 35 on both the "All Trades" sheet and the
# "base_strategy" sheet. Each sheet currently has a header in row 1 and
# trades in rows 2-34, so the new trade goes in row 35.

$wb = $excel.ActiveWorkbook

$sheetNames = @("All Trades", "base_strategy")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $row = 35

    $ws.Cells.Item($row, 1).Value = 34

    # Column B holds a plain "yyyy-mm-dd" text label (not a real date).
    # Force the cell to Text format first so Excel's automatic date
    # recognition doesn't silently turn the literal string into a date
    # serial number, then drop the format back to the sheet's normal
    # (unstyled) look - only the cell's text content should differ from
    # its neighbours, not its formatting.
    $ws.Cells.Item($row, 2).NumberFormat = "@"
    $ws.Cells.Item($row, 2).Value = "2026-02-16"
    $ws.Cells.Item($row, 2).Style = "Normal"

    $ws.Cells.Item($row, 3).Value = "22:55:29"
    $ws.Cells.Item($row, 4).Value = "base_strategy"
    $ws.Cells.Item($row, 5).Value = "DOWN"
    $ws.Cells.Item($row, 6).Value = 49.999998
    # Exit Price - trade is still OPEN, so this column has no value yet.
    $ws.Cells.Item($row, 7).Value = ""
    $ws.Cells.Item($row, 8).Value = "OPEN"
    $ws.Cells.Item($row, 9).Value = 0
    $ws.Cells.Item($row, 10).Value = 0
    $ws.Cells.Item($row, 11).Value = 100
    $ws.Cells.Item($row, 12).Value = 0
    $ws.Cells.Item($row, 13).Value = 0
    $ws.Cells.Item($row, 14).Value = 0.6
    $ws.Cells.Item($row, 15).Value = "Normal spread capture: 19600 bps"
    # Exit Reason - trade is still OPEN, so this column has no value yet.
    $ws.Cells.Item($row, 16).Value = ""
    $ws.Cells.Item($row, 17).Value = 0
}
